$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing rows down by one (row1->row2 ... row33->row34)
$ws.Rows("1:1").Insert()

# Update header cells: "Unit - activity1" -> "Period", "PJ" -> "2050"
# Apostrophe-prefix keeps these as text (matching original string type/shared string)
$ws.Range("B2").Value = "'Period"
$ws.Range("B3").Value = "'2050"

# Update the per-country H2 demand values for 2050 (column B, rows 4-34)
$ws.Range("B4").Value = 165.34299999999999
$ws.Range("B5").Value = 384.41
$ws.Range("B6").Value = 7.8729999999999993
$ws.Range("B7").Value = 37.160999999999994
$ws.Range("B8").Value = 7.7519999999999998
$ws.Range("B9").Value = 152.876
$ws.Range("B10").Value = 1078.385
$ws.Range("B11").Value = 82.462999999999994
$ws.Range("B12").Value = 16.077999999999999
$ws.Range("B13").Value = 131.37700000000001
$ws.Range("B14").Value = 1060.902
$ws.Range("B15").Value = 122.435
$ws.Range("B16").Value = 813.42399999999998
$ws.Range("B17").Value = 14.442
$ws.Range("B18").Value = 71.384
$ws.Range("B19").Value = 17.198
$ws.Range("B20").Value = 3.5739999999999998
$ws.Range("B21").Value = 403.58500000000004
$ws.Range("B22").Value = 4.3109999999999999
$ws.Range("B23").Value = 43.923999999999999
$ws.Range("B24").Value = 32.981999999999999
$ws.Range("B25").Value = 1.3879999999999999
$ws.Range("B26").Value = 617.89699999999993
$ws.Range("B27").Value = 30.905000000000001
$ws.Range("B28").Value = 431.07100000000003
$ws.Range("B29").Value = 26.777999999999999
$ws.Range("B30").Value = 59.68
$ws.Range("B31").Value = 128.23699999999999
$ws.Range("B32").Value = 34.485000000000007
$ws.Range("B33").Value = 101.143
$ws.Range("B34").Value = 606.94999999999993

# Narrow column B to fit the new (shorter) values/header
$ws.Columns("B:B").ColumnWidth = 4.65

# Restore the selected cell shown in the saved view
$ws.Range("D8").Select()
